$wb = $excel.ActiveWorkbook

$wsTimeLog = $wb.Worksheets.Item("time log")
$wsTasks   = $wb.Worksheets.Item("tasks")
$wsLinks   = $wb.Worksheets.Item("useful links")

# ---------------------------------------------------------------------
# "time log" sheet: log two new entries (rows 19 and 20)
# ---------------------------------------------------------------------
$wsTimeLog.Range("E19").Value = 43396
$wsTimeLog.Range("E19").NumberFormat = "d-mmm-yy"
$wsTimeLog.Range("F19").Value = "solder pins to kl25z for imu function"
$wsTimeLog.Range("I19").Value = 1.5
$wsTimeLog.Rows("19").RowHeight = 29

$wsTimeLog.Range("E20").Value = 43397
$wsTimeLog.Range("E20").NumberFormat = "d-mmm-yy"
$wsTimeLog.Range("F20").Value = "read data from imu on kl25z"
$wsTimeLog.Range("I20").Value = 2

# ---------------------------------------------------------------------
# "useful links" sheet: new reference row for the KL25Z platform page
# ---------------------------------------------------------------------
$wsLinks.Range("D11").Value = "https://os.mbed.com/platforms/KL25Z/"
$wsLinks.Range("E11").Value = "KL25Z info"

# ---------------------------------------------------------------------
# View/selection state for each sheet
# ---------------------------------------------------------------------
$wsTimeLog.Activate()
$winTimeLog = $excel.ActiveWindow
$winTimeLog.ScrollRow = 10
$winTimeLog.ScrollColumn = 1
$wsTimeLog.Range("I21").Select() | Out-Null

$wsTasks.Activate()
$winTasks = $excel.ActiveWindow
$winTasks.ScrollRow = 4
$winTasks.ScrollColumn = 1
$wsTasks.Range("D9").Select() | Out-Null

# "useful links" ends up the active/visible tab, matching activeTab="2"
$wsLinks.Activate()
$wsLinks.Range("E12").Select() | Out-Null
